# Update crypto price/volume data per latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.658.87'
$ws.Range('E2').Value = '  +3.60%  '
$ws.Range('D3').Value = '1.701.02'
$ws.Range('E3').Value = '  +2.57%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.56%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3947'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4041'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.544'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '55.40'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.94%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08821'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.305'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.46'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001336'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.648'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.38%  '
$ws.Range('D17').Value = '1.704.04'
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '101.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07071'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.932'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.28%  '
$ws.Range('D24').Value = '24.649.99'
$ws.Range('E24').Value = '  +3.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.995'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.339'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '159.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.251'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.12'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.679'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +17.77%  '
$ws.Range('E32').Value = '  -1.64%  '
$ws.Range('D33').Value = '1.886.37'
$ws.Range('E33').Value = '  +2.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.386'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +13.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08579'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.27'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2766'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.950'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.82'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02816'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09069'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.472'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7788'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7302'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.57'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.531'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.216'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.381'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +20.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.000'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.95'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.08046'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.81%  '
